$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) to the s_vals sheet, matching the existing
# header row's formatting (bold, centered, bordered) by copying the format
# from the adjacent header cell (G1), then set the data row value.

$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("H2").Value = 0

$excel.CutCopyMode = $false
